# Update database and shift read_price algorithm:
# Drop the oldest quarter (column D, "فصل دوم منتهی به 1399/06") and append the new
# quarter ("فصل چهارم منتهی به 1401/12") as the new last data column (M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest quarter column (D) - remaining columns shift left automatically,
# sliding each column's width/format one slot to the left as well.
$ws.Range("D1:D28").Delete() | Out-Null

# The freed-up column M needs the width that the new (year-end) column should carry.
$ws.Range("M1:M28").ColumnWidth = 30.17

# Give the freshly appended column M the same cell formatting as its left
# neighbour (column L) row by row, then fill in the new quarter's data.
$rowsToFill = 8, 9, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27

foreach ($r in $rowsToFill) {
    $ws.Range("L$r").Copy() | Out-Null
    $ws.Range("M$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}
$excel.CutCopyMode = $false

# New quarter header + publish date for the newly appended column.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-13 (9)"

# New quarter's financial figures.
$ws.Range("M11").Value = 1875266
$ws.Range("M12").Value = -1150155
$ws.Range("M13").Value = 725111
$ws.Range("M14").Value = -89402
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("M17").Value = 635709
$ws.Range("M18").Value = -11024
$ws.Range("M19").Value = 19149
$ws.Range("M20").Value = 643834
$ws.Range("M21").Value = 25561
$ws.Range("M22").Value = 669395
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 669395
$ws.Range("M25").Value = 3347
$ws.Range("M26").Value = 200000
$ws.Range("M27").Value = 3347
